$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 45882
$ws.Range("I4").Value = 45883

# Copy the date style from I4 onto I5 so the new date cell matches the
# formatting used by the rest of the date column.
$ws.Range("I4").Copy()
$ws.Range("I5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I5").Value = 45884
$ws.Range("J5").Value = 9.81
$ws.Range("K5").Value = 0.74
$ws.Range("L5").Value = 9.07
